$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cells I1 and J1 (copy style from existing header cell H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$excel.CutCopyMode = $false

# Data values for I2:J46
$iValues = @(7,6,7,7,7,6,8,8,10,9,8,8,8,8,8,8,8,7,8,8,8,8,8,8,8,8,9,7,8,6,8,6,7,7,7,7,6,5,7,7,4,4,4,5,4)
$jValues = @(8,7,7,7,7,7,9,8,10,9,8,8,8,8,8,8,8,8,8,8,8,8,8,8,8,8,9,8,9,6,8,7,7,7,7,7,6,5,7,7,5,4,4,5,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
